$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Replace "Ready for handoff" with "In Translation" wherever it occurs
$ws1.Range("E2:F4").Value = "In Translation"
$ws2.Range("C2:C4").Value = "In Translation"
$ws3.Range("C2:C4").Value = "In Translation"
